$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: dimension/measure classifications curated
$ws.Range("A2").Value = "iaest-measure:estructura-hogar"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "sdmx-dimension:refArea"

# Row 3: dim/medida swapped for the curated columns
$ws.Range("A3").Value = "medida"
$ws.Range("D3").Value = "dim"

# Row 4: concept/type column now carries URI mapping references
$ws.Range("A4").Value = "xsd:int"
$ws.Range("D4").Value = "URI-Municipio"
$ws.Range("F4").Value = "URI-Comunidad"

# Row 5 (mapping file references) is removed entirely
$ws.Range("A5:I5").Delete()
